$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "notes"
$ws.Range("H6").Value = "first good GB data 1699833600"
$ws.Range("H7").Value = "2 generators on 1 GB, not handled yet"

$ws.Columns.Item(8).ColumnWidth = 52.88671875

$ws.Range("H8").Select()
